$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 290; this shifts the existing rows 290-297
# down to 291-298 (old row 297 duplicates down to the new row 298),
# matching the diff's net effect of growing the sheet from R297 to R298.
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row 290 with the new weekly record.
$ws.Cells.Item(290, 1).Value = 2
$ws.Cells.Item(290, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(290, 3).Value = "Coquimbo"
$ws.Cells.Item(290, 4).Value = 45239
$ws.Cells.Item(290, 5).Value = 4
$ws.Cells.Item(290, 6).Value = 100112031
$ws.Cells.Item(290, 7).Value = "Poroto verde"
$ws.Cells.Item(290, 8).Value = "Magnum"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 500
$ws.Cells.Item(290, 11).Value = 45000
$ws.Cells.Item(290, 12).Value = 50000
$ws.Cells.Item(290, 13).Value = 47500
$ws.Cells.Item(290, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(290, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(290, 16).Value = 1900
$ws.Cells.Item(290, 17).Value = 25
$ws.Cells.Item(290, 18).Value = "Hortaliza"
